$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their literal text representation (avoid Excel
# auto-converting numeric-looking strings like "561.47" into floating point
# numbers, which would lose precision / formatting).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '65.651.50'
$ws.Range('E2').Value = '  +3.36%  '
$ws.Range('D3').Value = '3.409.54'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '561.47'
$ws.Range('E5').Value = '  +2.05%  '
$ws.Range('D6').Value = '175.24'
$ws.Range('E6').Value = '  +2.52%  '
$ws.Range('E7').Value = '  +2.85%  '
$ws.Range('D8').Value = '3.401.78'
$ws.Range('E8').Value = '  +2.32%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +14.00%  '
$ws.Range('D11').Value = '0.631'
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').Value = '55.11'
$ws.Range('E12').Value = '  +3.68%  '
$ws.Range('E13').Value = '  +6.70%  '
$ws.Range('D14').Value = '9.14'
$ws.Range('E14').Value = '  +2.81%  '
$ws.Range('D15').Value = '3.960.34'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = '18.34'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').Value = '3.416.82'
$ws.Range('E17').Value = '  +3.38%  '
$ws.Range('E18').Value = '  +1.54%  '
$ws.Range('D19').Value = '65.614.31'
$ws.Range('E19').Value = '  +3.52%  '
$ws.Range('D20').Value = '11.92'
$ws.Range('E20').Value = '  +2.60%  '
$ws.Range('D21').Value = '0.995'
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').Value = '469.20'
$ws.Range('E22').Value = '  +15.92%  '
$ws.Range('D23').Value = '5.06'
$ws.Range('E23').Value = '  +19.06%  '
$ws.Range('D24').Value = '4.14'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').Value = '86.35'
$ws.Range('E25').Value = '  +4.55%  '
$ws.Range('D26').Value = '13.63'
$ws.Range('E26').Value = '  +2.94%  '
$ws.Range('E27').Value = '  +3.03%  '
$ws.Range('D28').Value = '2.89'
$ws.Range('E28').Value = '  +6.33%  '
$ws.Range('D29').Value = '8.88'
$ws.Range('E29').Value = '  +3.64%  '
$ws.Range('D30').Value = '30.84'
$ws.Range('E30').Value = '  +6.32%  '
$ws.Range('D31').Value = '6.69'
$ws.Range('E31').Value = '  +4.71%  '
$ws.Range('D32').Value = '11.55'
$ws.Range('E32').Value = '  +2.22%  '
$ws.Range('D33').Value = '587.69'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '61.47'
$ws.Range('E34').Value = '  +7.34%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.109'
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  -3.29%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '3.50'
$ws.Range('E38').Value = '  +2.84%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '35.90'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').Value = '0.0₃0759'
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('D42').Value = '3.102.38'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = '2.86'
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('D45').Value = '0.0416'
$ws.Range('E45').Value = '  +3.82%  '
$ws.Range('D46').Value = '2.50'
$ws.Range('E46').Value = '  +2.94%  '
$ws.Range('D47').Value = '3.22'
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('D48').Value = '0.135'
$ws.Range('E48').Value = '  +5.86%  '
$ws.Range('D49').Value = '2.59'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').Value = '8.34'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('D51').Value = '136.18'
$ws.Range('E51').Value = '  +2.79%  '
